$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New journal entry - Day 3 of N8N learning (2025-09-04)
$ws.Range("A4").Value = 45904
$ws.Range("A4").NumberFormat = "DD/MM/YY"

$ws.Range("B4").Value = "GitHub Push Email Notifier"
$ws.Range("C4").Value = "Sends an email when a new push happens in my GitHub repo."
$ws.Range("D4").Value = "n8n_workflow_notifier.json"
